$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while preserving it as literal Text (matching the
# original inline-string cell type) and without leaving a permanent
# NumberFormat/style change behind. Many "Price" values look numeric
# (e.g. "1.00", "26.204.56", "142.80") but must round-trip byte-for-byte as
# text, so we briefly force Text format, assign the raw string via Value2
# (which skips Excel's "looks like a number" coercion heuristics less
# aggressively than .Value), then restore the cell's original Style object so
# no stray formatting diff is introduced.
function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Val
    )
    $cell = $ws.Range($Addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $Val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '26.204.56'
Set-TextValue 'E2' '  +1.66%  '
Set-TextValue 'D3' '1.644.27'
Set-TextValue 'E3' '  +0.54%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  -0.17%  '
Set-TextValue 'D5' '217.13'
Set-TextValue 'E6' '  +1.18%  '
Set-TextValue 'E7' '  -0.07%  '
Set-TextValue 'E8' '  +0.58%  '
Set-TextValue 'E9' '  +0.46%  '
Set-TextValue 'D10' '19.81'
Set-TextValue 'E10' '  +1.46%  '
Set-TextValue 'E11' '  +0.12%  '
Set-TextValue 'D12' '1.872.12'
Set-TextValue 'E12' '  +0.59%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '4.29'
Set-TextValue 'E13' '  +0.75%  '
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.641.44'
Set-TextValue 'E14' '  +0.35%  '
Set-TextValue 'D15' '0.545'
Set-TextValue 'E15' '  -3.11%  '
Set-TextValue 'D16' '0.0₃0765'
Set-TextValue 'E16' '  +0.17%  '
Set-TextValue 'D17' '63.27'
Set-TextValue 'E17' '  +0.08%  '
Set-TextValue 'D18' '26.211.28'
Set-TextValue 'E18' '  +1.56%  '
Set-TextValue 'E19' '  -0.13%  '
Set-TextValue 'B20' 'Uniswap'
Set-TextValue 'C20' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D20' '4.43'
Set-TextValue 'E20' '  -0.67%  '
Set-TextValue 'B21' 'BitcoinCash'
Set-TextValue 'C21' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D21' '195.14'
Set-TextValue 'E21' '  +1.44%  '
Set-TextValue 'E22' '  +0.82%  '
Set-TextValue 'E23' '  -0.49%  '
Set-TextValue 'E24' '  -3.07%  '
Set-TextValue 'E25' '  -0.11%  '
Set-TextValue 'D26' '142.80'
Set-TextValue 'E26' '  +0.69%  '
Set-TextValue 'E27' '  +0.89%  '
Set-TextValue 'E28' '  +0.57%  '
Set-TextValue 'D29' '15.63'
Set-TextValue 'E29' '  +0.90%  '
Set-TextValue 'E30' '  +0.90%  '
Set-TextValue 'D31' '0.0501'
Set-TextValue 'E31' '  +1.84%  '
Set-TextValue 'E32' '  +0.97%  '
Set-TextValue 'E33' '  +0.70%  '
Set-TextValue 'E34' '  +1.66%  '
Set-TextValue 'D35' '2.42'
Set-TextValue 'D36' '0.912'
Set-TextValue 'E36' '  +1.03%  '
Set-TextValue 'D37' '1.134.36'
Set-TextValue 'E37' '  +0.28%  '
Set-TextValue 'D38' '0.553'
Set-TextValue 'E38' '  +1.64%  '
Set-TextValue 'E39' '  -0.96%  '
Set-TextValue 'E40' '  +1.22%  '
Set-TextValue 'D41' '1.00'
Set-TextValue 'E41' '  -0.17%  '
Set-TextValue 'D42' '100.56'
Set-TextValue 'E42' '  -0.15%  '
Set-TextValue 'D43' '5.49'
Set-TextValue 'E43' '  -1.21%  '
Set-TextValue 'D44' '0.798'
Set-TextValue 'E44' '  -0.46%  '
Set-TextValue 'D45' '1.781.61'
Set-TextValue 'E45' '  +0.63%  '
Set-TextValue 'E46' '  -0.64%  '
Set-TextValue 'D47' '56.94'
Set-TextValue 'E47' '  +2.93%  '
Set-TextValue 'B48' 'Cronos'
Set-TextValue 'C48' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D48' '0.0517'
Set-TextValue 'E48' '  +2.98%  '
Set-TextValue 'B49' 'RenderToken'
Set-TextValue 'C49' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D49' '1.47'
Set-TextValue 'E49' '  +3.24%  '
Set-TextValue 'B50' 'Mantle'
Set-TextValue 'C50' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D50' '0.417'
Set-TextValue 'E50' '  +0.16%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.69'
Set-TextValue 'E51' '  +3.01%  '
